$d = $word.ActiveDocument

# --- Change 1: paragraph 2 ---
# "...can be tricky, but I will show you a trick to make this easier!"
# -> "...can be tricky, but there's a special trick to make it easier!"
$d.Content.Find.Execute(
    "I will show you a trick to make this easier",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "there’s a special trick to make it easier", 2) | Out-Null

# --- Change 2: paragraph 3 ---
# "...math concept that is used for this special trick."
# -> "...math concept that we will be using for this trick."
$d.Content.Find.Execute(
    "is used for this special trick.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "we will be using for this trick.", 2) | Out-Null

# --- Change 3: paragraph 5 ---
# "...we can split a large number up into smaller pieces, compute each one, and add it all back together."
# -> "...we can split up the large number of an equation into smaller pieces."
$d.Content.Find.Execute(
    "we can split a large number up into smaller pieces, compute each one, and add it all back together.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "we can split up the large number of an equation into smaller pieces.", 2) | Out-Null

# --- Change 4: paragraph 10 ---
# "And finally, we add those two areas..." -> "And finally, we add these two areas..."
$d.Content.Find.Execute(
    "And finally, we add those two areas together",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "And finally, we add these two areas together", 2) | Out-Null

# --- Change 5: move "[display drag instruction and visual indicator]" paragraph
#     from after "Let's give it a try..." to before it ---
$tryRange = $d.Content
$tryRange.Find.Execute("Let’s give it a try", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tryPara = $tryRange.Paragraphs(1)
$tryIndex = $tryPara.Index

$tryPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs($tryIndex)
$newPara.Range.InsertBefore("[display drag instruction and visual indicator]")

$dupPara = $d.Paragraphs($tryIndex + 2)
$dupPara.Range.Delete()

# --- Change 6: split the "Finally, we add the products..." paragraph into two ---
$splitRange = $d.Content
$splitRange.Find.Execute("final product!", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endPos = $splitRange.End
$insertPoint = $d.Range($endPos, $endPos)
$insertPoint.InsertParagraphAfter()
$d.Range($endPos + 1, $endPos + 2).Delete()
